$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, even when the text looks like
# a number (e.g. "228.61" or "18.40"). Excel's normal Value-assignment
# auto-detects numeric-looking strings and stores them as numbers, which
# would both change the stored type and silently drop formatting like
# trailing zeros. Forcing a text NumberFormat before the assignment keeps
# the value a genuine string; resetting the style back to Normal afterwards
# avoids leaving a stray "text" style applied to the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "38.796.92"
$ws.Range("E2").Value = "  +2.91%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.097.70"
$ws.Range("E3").Value = "  +2.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "228.61"
$ws.Range("E5").Value = "  +0.66%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.14%  "

# Row 7 - Solana
Set-TextValue "D7" "60.40"
$ws.Range("E7").Value = "  +1.57%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.385"
$ws.Range("E9").Value = "  +2.15%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0836"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.28%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.408.32"
$ws.Range("E12").Value = "  +3.09%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +3.82%  "

# Row 14 - Avalanche
Set-TextValue "D14" "22.12"
$ws.Range("E14").Value = "  +5.31%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.799"
$ws.Range("E15").Value = "  +3.79%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.47"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.102.90"
$ws.Range("E17").Value = "  +3.24%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "38.726.18"
$ws.Range("E18").Value = "  +2.91%  "

# Row 19 - Litecoin
Set-TextValue "D19" "72.04"
$ws.Range("E19").Value = "  +3.81%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.06"
$ws.Range("E20").Value = "  +2.65%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.67%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "226.33"
$ws.Range("E22").Value = "  +1.34%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.42"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.61%  "

# Row 26 - was Cosmos, now Monero
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "170.45"
$ws.Range("E26").Value = "  +1.34%  "

# Row 27 - was Monero, now Cosmos
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "9.51"
$ws.Range("E27").Value = "  +1.54%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +5.94%  "

# Row 29 - ImmutableX
Set-TextValue "D29" "1.40"
$ws.Range("E29").Value = "  +9.79%  "

# Row 30 - EthereumClassic
$ws.Range("E30").Value = "  +2.22%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +0.76%  "

# Row 32 - WEMIXToken
Set-TextValue "D32" "2.35"
$ws.Range("E32").Value = "  +5.35%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "4.76"
$ws.Range("E33").Value = "  +6.38%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +2.64%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0614"
$ws.Range("E35").Value = "  +1.44%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +3.16%  "

# Row 37 - THORChain
$ws.Range("E37").Value = "  +1.47%  "

# Row 38 - RenderToken
Set-TextValue "D38" "3.51"
$ws.Range("E38").Value = "  +2.41%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  +0.08%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "18.40"
$ws.Range("E40").Value = "  +2.14%  "

# Row 41 - Aave
Set-TextValue "D41" "101.63"
$ws.Range("E41").Value = "  +4.59%  "

# Row 42 - Maker (only price changes)
Set-TextValue "D42" "1.541.38"

# Row 43 - VeChain
$ws.Range("E43").Value = "  +3.51%  "

# Row 44 - Cronos
$ws.Range("E44").Value = "  +2.50%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  -0.73%  "

# Row 46 - FraxShare
Set-TextValue "D46" "7.69"
$ws.Range("E46").Value = "  +9.72%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +1.10%  "

# Row 48 - FTXToken
Set-TextValue "D48" "4.10"
$ws.Range("E48").Value = "  -5.53%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +3.23%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.294.85"
$ws.Range("E51").Value = "  +3.14%  "
